$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28, shifting existing rows 28-79
# down to 29-80 (preserving their data/formatting untouched).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44581
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 100112026
$ws.Range("G28").Value = "Haba"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 25000
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 1000
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
